$d = $word.ActiveDocument

$replacements = @(
    @("639×9=5751", "192×3=576"),
    @("283×6=1698", "902×8=7216"),
    @("616×5=3080", "320×8=2560"),
    @("124×4=496", "433×8=3464"),
    @("538×8=4304", "842×5=4210"),
    @("766×6=4596", "772×2=1544"),
    @("724×4=2896", "772×3=2316"),
    @("822×5=4110", "801×9=7209"),
    @("151×5=755", "933×5=4665"),
    @("714×3=2142", "640×9=5760"),
    @("736×5=3680", "514×7=3598"),
    @("761×7=5327", "401×8=3208"),
    @("450×5=2250", "575×2=1150"),
    @("744×7=5208", "206×5=1030"),
    @("687×2=1374", "149×8=1192"),
    @("885×8=7080", "772×4=3088"),
    @("542×8=4336", "396×6=2376"),
    @("191×4=764", "855×8=6840"),
    @("113×7=791", "151×9=1359"),
    @("106×5=530", "389×3=1167"),
    @("850×8=6800", "156×5=780"),
    @("286×2=572", "196×7=1372"),
    @("847×9=7623", "455×7=3185"),
    @("811×5=4055", "203×6=1218"),
    @("508×2=1016", "123×2=246")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
